$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9593005776405334
$ws.Range("B1").Value = 4.524675369262695
$ws.Range("C1").Value = 4.172035217285156
$ws.Range("D1").Value = 2.453616142272949
$ws.Range("E1").Value = 2.04694390296936
